$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is brand new - copy the formatting (styles) from row 4 (columns A:J only,
# matching the target which has no K5 cell) before filling in values.
$ws.Range("A4:J4").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 5: Crumpet GEF / Crumpet exporter ---
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# --- Row 6: Scone GEF / Scone exporter (row already existed, just empty) ---
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Update the selection to reflect the newly entered rows
$ws.Range("A5:J6").Select() | Out-Null
